# Adds 60 new ranking rows (98-157) to Sheet1, updates C2's score,
# and selects the first empty row below the appended data,
# reproducing the uploaded workbook's new state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated point total for the first ranked entry (row 2, column C).
$ws.Range("C2").Value = 18.010000000000002

# New name/points pairs appended to the ranking table (rows 98-157).
$newRows = @(
    @("הגר אגמון", 1),
    @("אורי שטרנברג", 1),
    @("ירון גלפנד", 1),
    @("יובל סטרוזר", 1),
    @("איתי הראל", 1),
    @("אורי שטרנברג", 6),
    @("יובל סטרוזר", 6),
    @("אביב ואסקז", 1),
    @("עדן ורד מרי", 1),
    @("תאיו ורד", 1),
    @("איתי הראל", 1),
    @("דפנה ברגשטיין", 1),
    @("הילה שולויס", 1),
    @("יהלי גודר", 1),
    @("יער אלביר", 1),
    @("ליאם דיין ", 1),
    @("מעיין סטרוזר", 1),
    @("עמית גורוביץ", 1),
    @("תאיו ורד", 6),
    @("ליאם דיין ", 6),
    @("אביב ואסקז", 1),
    @("הגר אגמון", 1),
    @("אורי שטרנברג", 1),
    @("ירון גלפנד", 1),
    @("ליאם דיין ", 1),
    @("ליהי בראל", 1),
    @("יהלי דוייב", 1),
    @("תומר ששון", 1),
    @("יער אלביר", 1),
    @("ירון גלפנד", 6),
    @("ירון גלפנד", 6),
    @("לידור אלשטיין", 1),
    @("דן פימה", 1),
    @("תאיו ורד", 1),
    @("יובל סטרוזר", 1),
    @("מעיין סטרוזר", 1),
    @("תומר ששון", 1),
    @("יולי יערי תליו", 1),
    @("איתי הראל", 1),
    @("דפנה ברגשטיין", 1),
    @("תומר ששון", 6),
    @("לידור אלשטיין", 6),
    @("לידור אלשטיין", 1),
    @("מעיין סטרוזר", 1),
    @("גלי זליג", 1),
    @("יובל סטרוזר", 1),
    @("גלי זליג", 6),
    @("גלי זליג", 6),
    @("אביב ואסקז", 1),
    @("עדן ורד מרי", 1),
    @("ליהי בראל", 1),
    @("תאיו ורד", 1),
    @("יער אלביר", 1),
    @("תומר ששון", 1),
    @("דפנה ברגשטיין", 1),
    @("ליאם דיין ", 1),
    @("יהלי דוייב", 1),
    @("איתי הראל", 1),
    @("ליהי בראל", 6),
    @("אביב ואסקז", 6)
)

$startRow = 98
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $name = $newRows[$i][0]
    $points = $newRows[$i][1]
    $ws.Cells.Item($rowNum, 1).Value = $name
    $ws.Cells.Item($rowNum, 2).Value = $points
}

# Scroll/select to match the author's final cursor position (first blank row).
$ws.Range("A158").Select()
